$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144, shifting existing rows 144:173 down to 145:174.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new "Sweet Heart" record.
$ws.Cells.Item(144, 1).Value = 7
$ws.Cells.Item(144, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(144, 3).Value = "Ñuble"
$ws.Cells.Item(144, 4).Value = 44951
$ws.Cells.Item(144, 5).Value = 16
$ws.Cells.Item(144, 6).Value = "Fruta"
$ws.Cells.Item(144, 7).Value = 100103
$ws.Cells.Item(144, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(144, 9).Value = 100103001
$ws.Cells.Item(144, 10).Value = "Cereza"
$ws.Cells.Item(144, 11).Value = "Sweet Heart"
$ws.Cells.Item(144, 12).Value = "Primera"
$ws.Cells.Item(144, 13).Value = 50
$ws.Cells.Item(144, 14).Value = 8000
$ws.Cells.Item(144, 15).Value = 8000
$ws.Cells.Item(144, 16).Value = 8000
$ws.Cells.Item(144, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(144, 18).Value = "Quillón"
$ws.Cells.Item(144, 19).Value = 800
$ws.Cells.Item(144, 20).Value = 10
